# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# values on the zh-cn and de-de sheets to reflect the newly generated report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-31 07:39:58"
$wsZhCn.Range("H2").Value = "2016-03-31 07:40:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-31 07:40:10"
$wsDeDe.Range("H2").Value = "2016-03-31 07:41:09"
